$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows of data (Common Name, TOOP Document Type Identifier, Since, Deprecated?)
# matching the "Added new WP3 IDs" commit.
$rows = @(
    @("Ship Certificate",           "urn:eu:toop:ns:dataexchange-1p40::Request##urn:eu.toop.request.shipcertificate-list::1.40"),
    @("Ship Certificate",           "urn:eu:toop:ns:dataexchange-1p40::Response##urn:eu.toop.response.shipcertificate-list::1.40"),
    @("Ship Certificate",           "urn:eu:toop:ns:dataexchange-1p40::Request##urn:eu.toop.request.shipcertificate::1.40"),
    @("Ship Certificate",           "urn:eu:toop:ns:dataexchange-1p40::Response##urn:eu.toop.response.shipcertificate::1.40"),
    @("Crew Certificate",           "urn:eu:toop:ns:dataexchange-1p40::Request##urn:eu.toop.request.crewcertificate-list::1.40"),
    @("Crew Certificate",           "urn:eu:toop:ns:dataexchange-1p40::Response##urn:eu.toop.response.crewcertificate-list::1.40"),
    @("Crew Certificate",           "urn:eu:toop:ns:dataexchange-1p40::Request##urn:eu.toop.request.crewcertificate::1.40"),
    @("Crew Certificate",           "urn:eu:toop:ns:dataexchange-1p40::Response##urn:eu.toop.response.crewcertificate::1.40"),
    @("Registered Organization",    "urn:eu:toop:ns:dataexchange-1p40::Request##urn:eu.toop.request.registeredorganization-list::1.40"),
    @("Registered Organization",    "urn:eu:toop:ns:dataexchange-1p40::Response##urn:eu.toop.response.registeredorganization-list::1.40"),
    @("Registered Organization",    "urn:eu:toop:ns:dataexchange-1p40::Request##urn:eu.toop.request.registeredorganization::1.40"),
    @("Registered Organization",    "urn:eu:toop:ns:dataexchange-1p40::Response##urn:eu.toop.response.registeredorganization::1.40"),
    @("Evidence",                   "urn:eu:toop:ns:dataexchange-1p40::Request##urn:eu.toop.request.evidence-list::1.40"),
    @("Evidence",                   "urn:eu:toop:ns:dataexchange-1p40::Response##urn:eu.toop.response.evidence-list::1.40"),
    @("Evidence",                   "urn:eu:toop:ns:dataexchange-1p40::Request##urn:eu.toop.request.evidence::1.40"),
    @("Evidence",                   "urn:eu:toop:ns:dataexchange-1p40::Response##urn:eu.toop.response.evidence::1.40")
)

$r = 6
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = 2
    $ws.Cells.Item($r, 4).Formula = "=FALSE"
    $r = $r + 1
}

# Wrap text on the identifier / common name columns (A & B) across the whole used range,
# matching the new column widths/behaviour in the refreshed sheet.
$ws.Range("A1:B21").WrapText = $true

# The "Registered Organization" / "Evidence" identifier cells keep the original
# left-aligned look; the new Ship/Crew Certificate rows do not.
$ws.Range("B2:B5").HorizontalAlignment = -4131
$ws.Range("B14:B21").HorizontalAlignment = -4131

$ws.Columns("A").ColumnWidth = 22.166666666666668
$ws.Columns("B").ColumnWidth = 93.16666666666667

# Move the active selection to A6, as captured in the saved workbook view.
$ws.Range("A6").Select() | Out-Null
